# Add a new "2022" column (P) to the table, mirroring the style of the
# existing "2021" column (O), and update the sheet's selection to Q4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header value for 2022 in row 4 (year header row)
$ws.Range("P4").Value = 2022

# New data values for the 2022 column, row by row
$ws.Range("P5").Value = 11.4
$ws.Range("P6").Value = 12.6
$ws.Range("P7").Value = 9.8
$ws.Range("P8").Value = 11.4
$ws.Range("P9").Value = 5.4
$ws.Range("P10").Value = 4.7
$ws.Range("P11").Value = 3.4
$ws.Range("P12").Value = 17.7
$ws.Range("P13").Value = 20.5
$ws.Range("P14").Value = 8.4
$ws.Range("P16").Value = 12.9
$ws.Range("P17").Value = 10.5

# Copy the formatting (number format, borders, font, fill, alignment) from
# the corresponding O-column cell into the new P-column cell so the new
# column visually matches its neighbour (O uses style index 6/28/27/30/34
# per row; copying from O keeps that mapping automatically).
$rowsToFormat = 4,5,6,7,8,9,10,11,12,13,14,16,17
foreach ($r in $rowsToFormat) {
    $src = $ws.Range("O$r")
    $dst = $ws.Range("P$r")
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

# Update the selection shown when the workbook is reopened
$ws.Range("Q4").Select()

# Update the worksheet's used-range dimension now includes column P
$ws.Range("A1:P21")
